# Insert a new data row at row 335 (pushing the existing rows 335-413 down
# to 336-414) and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(335).Insert()

$ws.Cells.Item(335, 1).Value = 9
$ws.Cells.Item(335, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(335, 3).Value = 'Metropolitana'
$ws.Cells.Item(335, 4).Value = 44798
$ws.Cells.Item(335, 5).Value = 13
$ws.Cells.Item(335, 6).Value = 100112032
$ws.Cells.Item(335, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(335, 8).Value = 'Sin especificar'
$ws.Cells.Item(335, 9).Value = 'Primera'
$ws.Cells.Item(335, 10).Value = 375
$ws.Cells.Item(335, 11).Value = 21000
$ws.Cells.Item(335, 12).Value = 24000
$ws.Cells.Item(335, 13).Value = 22573
$ws.Cells.Item(335, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(335, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(335, 16).Value = 451
$ws.Cells.Item(335, 17).Value = 50
$ws.Cells.Item(335, 18).Value = 'Hortaliza'
